# Update cryptocurrency price and 1h volume% values to the latest
# refresh pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.527.71"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.861.60"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Formula = "=""311.97"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Formula = "=""1.011"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Formula = "=""0.4777"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Formula = "=""0.3802"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("D9").Formula = "=""0.07334"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Formula = "=""0.9308"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Formula = "=""20.73"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").Formula = "=""0.07780"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "1.892.08"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Formula = "=""5.441"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Formula = "=""6.567"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Formula = "=""90.23"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Formula = "=""0.000008817"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "27.595.82"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Formula = "=""14.66"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Formula = "=""5.096"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Formula = "=""10.71"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Formula = "=""1.937"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Formula = "=""155.75"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Formula = "=""18.49"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Formula = "=""2.013"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Formula = "=""115.40"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Formula = "=""4.952"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Formula = "=""0.08864"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Formula = "=""3.328"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").Formula = "=""1.204"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").Formula = "=""0.7529"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").Formula = "=""4.583"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Formula = "=""2.689"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Formula = "=""1.122"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Formula = "=""0.02036"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +3.56%  "
$ws.Range("D38").Formula = "=""0.5639"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("D39").Formula = "=""0.05331"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").Formula = "=""2.980"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Formula = "=""7.018"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Formula = "=""0.1524"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Formula = "=""8.486"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").Formula = "=""10.72"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = "=""0.4869"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").Formula = "=""1.012"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Formula = "=""104.50"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").Formula = "=""67.43"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Formula = "=""0.9097"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +1.99%  "

$excel.CutCopyMode = $false
